$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: only B2 changes
$ws.Range("B2").Value = 2359532929403399

# Row 3 - RandomForestRegressor: B3, C3, D3 change
$ws.Range("B3").Value = 0.04685285517212339
$ws.Range("C3").Value = 0.04497891964514411
$ws.Range("D3").Value = 31884639615374.34

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor, B4, C4, D4 change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.04531815599783307
$ws.Range("C4").Value = 0.04321401098607317
$ws.Range("D4").Value = 0.07403079416724975

# Row 5 - AdaBoostRegressor -> MLPRegressor, B5, C5, D5 change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 57485378475428.41
$ws.Range("C5").Value = 12137759413514.89
$ws.Range("D5").Value = 28950205440960.89
